$d = $word.ActiveDocument

$replacements = @(
    @{old = "119÷3="; new = "333÷2="},
    @{old = "385÷3="; new = "946÷8="},
    @{old = "381÷8="; new = "288÷3="},
    @{old = "837÷4="; new = "205÷5="},
    @{old = "475÷8="; new = "490÷9="},
    @{old = "731÷6="; new = "313÷5="},
    @{old = "529÷8="; new = "846÷3="},
    @{old = "887÷9="; new = "127÷5="},
    @{old = "208÷2="; new = "913÷5="},
    @{old = "562÷5="; new = "124÷3="},
    @{old = "219÷7="; new = "623÷6="},
    @{old = "749÷4="; new = "571÷5="},
    @{old = "962÷7="; new = "642÷6="},
    @{old = "544÷6="; new = "779÷5="},
    @{old = "983÷8="; new = "506÷3="},
    @{old = "478÷7="; new = "647÷6="},
    @{old = "642÷4="; new = "261÷3="},
    @{old = "782÷8="; new = "477÷9="},
    @{old = "419÷4="; new = "462÷5="},
    @{old = "896÷4="; new = "464÷7="},
    @{old = "428÷3="; new = "280÷4="},
    @{old = "259÷3="; new = "168÷4="},
    @{old = "976÷6="; new = "177÷4="},
    @{old = "250÷9="; new = "868÷2="},
    @{old = "359÷2="; new = "449÷3="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
